# "Generate Report for Handoff"
#
# A fresh handoff run completed for the four "Ready for handoff" source
# files (2e17fcf0..., 7683bf05..., d7cb3f05..., ded995b5...). Refresh the
# generated-handoff timestamps and bump those rows' priority from "low" to
# "ht" across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G4:G7 -- "Latest HO Xliff Generate Date" for the 4 handed-off rows
$wsOverview.Range("G4:G7").Value = "2016-08-19 22:38:40"

# zh-cn sheet: Priority low -> ht, and Latest Handoff Datetime refreshed
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-08-19 22:38:36"

# de-de sheet: Priority low -> ht
# (de-de!H4:H7 "Latest Handoff Datetime" shares the same string as
#  Overview!G4:G7, so it is already refreshed to 22:38:40 above)
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-08-19 22:38:40"

Write-Host "Generate Report for Handoff: applied"
